# Replace the Ptn-Ptprb sending/target cluster combinations with the
# freshly recomputed TPM-based NATMI edge statistics (3x3 grid: ECs/FAPs/MuSCs).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Ptn"
$ws.Range("C2").Value = "Ptprb"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.06520933333333333
$ws.Range("H2").Value = 0.195628
$ws.Range("I2").Value = 0.007116089623850998
$ws.Range("J2").Value = 0.007116089623850999
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 154.017718
$ws.Range("N2").Value = 462.0531539999999
$ws.Range("O2").Value = 0.9743120958630357
$ws.Range("P2").Value = 0.9743120958630357
$ws.Range("Q2").Value = 10.04339271230133
$ws.Range("R2").Value = 90.390534410712
$ws.Range("S2").Value = 0.006933292195763467
$ws.Range("T2").Value = 0.006933292195763468

# Row 3: ECs -> FAPs
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Ptn"
$ws.Range("C3").Value = "Ptprb"
$ws.Range("D3").Value = "FAPs"
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.06520933333333333
$ws.Range("H3").Value = 0.195628
$ws.Range("I3").Value = 0.007116089623850998
$ws.Range("J3").Value = 0.007116089623850999
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 2.253965666666666
$ws.Range("N3").Value = 6.761896999999999
$ws.Range("O3").Value = 0.0142585284421194
$ws.Range("P3").Value = 0.0142585284421194
$ws.Range("Q3").Value = 0.1469795984795555
$ws.Range("R3").Value = 1.322816386316
$ws.Range("S3").Value = 0.0001014649662983502
$ws.Range("T3").Value = 0.0001014649662983502

# Row 4: ECs -> MuSCs
$ws.Range("A4").Value = "ECs"
$ws.Range("B4").Value = "Ptn"
$ws.Range("C4").Value = "Ptprb"
$ws.Range("D4").Value = "MuSCs"
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.06520933333333333
$ws.Range("H4").Value = 0.195628
$ws.Range("I4").Value = 0.007116089623850998
$ws.Range("J4").Value = 0.007116089623850999
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 1.806737666666667
$ws.Range("N4").Value = 5.420213
$ws.Range("O4").Value = 0.011429375694845
$ws.Range("P4").Value = 0.011429375694845
$ws.Range("Q4").Value = 0.1178161587515555
$ws.Range("R4").Value = 1.060345428764
$ws.Range("S4").Value = [double]"8.133246178918128e-05"
$ws.Range("T4").Value = [double]"8.133246178918128e-05"

# Row 5: FAPs -> ECs
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Ptn"
$ws.Range("C5").Value = "Ptprb"
$ws.Range("D5").Value = "ECs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 3.270036666666666
$ws.Range("H5").Value = 9.81011
$ws.Range("I5").Value = 0.3568488252184601
$ws.Range("J5").Value = 0.3568488252184602
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 154.017718
$ws.Range("N5").Value = 462.0531539999999
$ws.Range("O5").Value = 0.9743120958630357
$ws.Range("P5").Value = 0.9743120958630357
$ws.Range("Q5").Value = 503.6435851763266
$ws.Range("R5").Value = 4532.792266586939
$ws.Range("S5").Value = 0.3476821268048599
$ws.Range("T5").Value = 0.34768212680486

# Row 6: FAPs -> FAPs
$ws.Range("A6").Value = "FAPs"
$ws.Range("B6").Value = "Ptn"
$ws.Range("C6").Value = "Ptprb"
$ws.Range("D6").Value = "FAPs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 3.270036666666666
$ws.Range("H6").Value = 9.81011
$ws.Range("I6").Value = 0.3568488252184601
$ws.Range("J6").Value = 0.3568488252184602
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 2.253965666666666
$ws.Range("N6").Value = 6.761896999999999
$ws.Range("O6").Value = 0.0142585284421194
$ws.Range("P6").Value = 0.0142585284421194
$ws.Range("Q6").Value = 7.370550375407777
$ws.Range("R6").Value = 66.33495337867
$ws.Range("S6").Value = 0.005088139123914306
$ws.Range("T6").Value = 0.005088139123914308

# Row 7: FAPs -> MuSCs
$ws.Range("A7").Value = "FAPs"
$ws.Range("B7").Value = "Ptn"
$ws.Range("C7").Value = "Ptprb"
$ws.Range("D7").Value = "MuSCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 3.270036666666666
$ws.Range("H7").Value = 9.81011
$ws.Range("I7").Value = 0.3568488252184601
$ws.Range("J7").Value = 0.3568488252184602
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 1.806737666666667
$ws.Range("N7").Value = 5.420213
$ws.Range("O7").Value = 0.011429375694845
$ws.Range("P7").Value = 0.011429375694845
$ws.Range("Q7").Value = 5.908098417047777
$ws.Range("R7").Value = 53.17288575342999
$ws.Range("S7").Value = 0.004078559289685858
$ws.Range("T7").Value = 0.004078559289685858

# Row 8: MuSCs -> ECs
$ws.Range("A8").Value = "MuSCs"
$ws.Range("B8").Value = "Ptn"
$ws.Range("C8").Value = "Ptprb"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 5.828401
$ws.Range("H8").Value = 17.485203
$ws.Range("I8").Value = 0.6360350851576888
$ws.Range("J8").Value = 0.6360350851576889
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 154.017718
$ws.Range("N8").Value = 462.0531539999999
$ws.Range("O8").Value = 0.9743120958630357
$ws.Range("P8").Value = 0.9743120958630357
$ws.Range("Q8").Value = 897.6770216089179
$ws.Range("R8").Value = 8079.093194480262
$ws.Range("S8").Value = 0.6196966768624121
$ws.Range("T8").Value = 0.6196966768624123

# Row 9: MuSCs -> FAPs
$ws.Range("A9").Value = "MuSCs"
$ws.Range("B9").Value = "Ptn"
$ws.Range("C9").Value = "Ptprb"
$ws.Range("D9").Value = "FAPs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 5.828401
$ws.Range("H9").Value = 17.485203
$ws.Range("I9").Value = 0.6360350851576888
$ws.Range("J9").Value = 0.6360350851576889
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 2.253965666666666
$ws.Range("N9").Value = 6.761896999999999
$ws.Range("O9").Value = 0.0142585284421194
$ws.Range("P9").Value = 0.0142585284421194
$ws.Range("Q9").Value = 13.13701574556567
$ws.Range("R9").Value = 118.233141710091
$ws.Range("S9").Value = 0.009068924351906737
$ws.Range("T9").Value = 0.009068924351906739

# Row 10: MuSCs -> MuSCs
$ws.Range("A10").Value = "MuSCs"
$ws.Range("B10").Value = "Ptn"
$ws.Range("C10").Value = "Ptprb"
$ws.Range("D10").Value = "MuSCs"
$ws.Range("E10").Value = 3
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = 5.828401
$ws.Range("H10").Value = 17.485203
$ws.Range("I10").Value = 0.6360350851576888
$ws.Range("J10").Value = 0.6360350851576889
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.806737666666667
$ws.Range("N10").Value = 5.420213
$ws.Range("O10").Value = 0.011429375694845
$ws.Range("P10").Value = 0.011429375694845
$ws.Range("Q10").Value = 10.53039162313767
$ws.Range("R10").Value = 94.773524608239
$ws.Range("S10").Value = 0.007269483943369957
$ws.Range("T10").Value = 0.007269483943369957
